# Applies the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store these as literal text (quote-prefix),
# matching the source workbook's string cells -- without it, values such as
# '0.999' or '19.00' would be silently reinterpreted as numbers and lose
# their exact formatting (trailing zeros, multi-dot price groupings, etc.).
$q = "'"

# Row 2
$ws.Range('D2').Value = $q + '42.885.21'
$ws.Range('E2').Value = $q + '  -0.29%  '

# Row 3
$ws.Range('D3').Value = $q + '2.554.60'
$ws.Range('E3').Value = $q + '  +0.19%  '

# Row 4
$ws.Range('D4').Value = $q + '0.999'
$ws.Range('E4').Value = $q + '  -0.15%  '

# Row 5
$ws.Range('D5').Value = $q + '304.83'
$ws.Range('E5').Value = $q + '  +2.17%  '

# Row 6
$ws.Range('D6').Value = $q + '98.76'
$ws.Range('E6').Value = $q + '  +6.90%  '

# Row 7
$ws.Range('E7').Value = $q + '  +0.15%  '

# Row 8
$ws.Range('E8').Value = $q + '  +0.05%  '

# Row 9
$ws.Range('D9').Value = $q + '0.549'
$ws.Range('E9').Value = $q + '  -0.36%  '

# Row 10
$ws.Range('D10').Value = $q + '37.21'
$ws.Range('E10').Value = $q + '  +3.07%  '

# Row 11
$ws.Range('D11').Value = $q + '0.0810'
$ws.Range('E11').Value = $q + '  +0.09%  '

# Row 12
$ws.Range('B12').Value = $q + 'TRON'
$ws.Range('C12').Value = $q + 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = $q + '0.116'
$ws.Range('E12').Value = $q + '  +7.33%  '

# Row 13
$ws.Range('B13').Value = $q + 'Polkadot'
$ws.Range('C13').Value = $q + 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = $q + '7.64'
$ws.Range('E13').Value = $q + '  -0.48%  '

# Row 14
$ws.Range('D14').Value = $q + '2.564.08'
$ws.Range('E14').Value = $q + '  +0.48%  '

# Row 15
$ws.Range('D15').Value = $q + '15.01'
$ws.Range('E15').Value = $q + '  +5.85%  '

# Row 16
$ws.Range('D16').Value = $q + '0.881'
$ws.Range('E16').Value = $q + '  +1.09%  '

# Row 17
$ws.Range('D17').Value = $q + '42.952.12'
$ws.Range('E17').Value = $q + '  -0.22%  '

# Row 18
$ws.Range('D18').Value = $q + '13.81'
$ws.Range('E18').Value = $q + '  +7.60%  '

# Row 19
$ws.Range('D19').Value = $q + '0.0₃0989'
$ws.Range('E19').Value = $q + '  +0.99%  '

# Row 20
$ws.Range('D20').Value = $q + '6.65'
$ws.Range('E20').Value = $q + '  -0.03%  '

# Row 21
$ws.Range('D21').Value = $q + '71.79'
$ws.Range('E21').Value = $q + '  +0.04%  '

# Row 22
$ws.Range('D22').Value = $q + '255.67'
$ws.Range('E22').Value = $q + '  -1.79%  '

# Row 23
$ws.Range('D23').Value = $q + '2.98'
$ws.Range('E23').Value = $q + '  +2.46%  '

# Row 24
$ws.Range('D24').Value = $q + '2.10'
$ws.Range('E24').Value = $q + '  -2.09%  '

# Row 25
$ws.Range('D25').Value = $q + '27.88'
$ws.Range('E25').Value = $q + '  -5.61%  '

# Row 26
$ws.Range('E26').Value = $q + '  -0.12%  '

# Row 27
$ws.Range('D27').Value = $q + '10.16'
$ws.Range('E27').Value = $q + '  +1.03%  '

# Row 28
$ws.Range('D28').Value = $q + '38.22'
$ws.Range('E28').Value = $q + '  +3.97%  '

# Row 29
$ws.Range('E29').Value = $q + '  -1.46%  '

# Row 30
$ws.Range('E30').Value = $q + '  +0.76%  '

# Row 31
$ws.Range('D31').Value = $q + '158.51'
$ws.Range('E31').Value = $q + '  +2.57%  '

# Row 32
$ws.Range('E32').Value = $q + '  +0.07%  '

# Row 33
$ws.Range('E33').Value = $q + '  +0.92%  '

# Row 34
$ws.Range('D34').Value = $q + '0.0810'
$ws.Range('E34').Value = $q + '  +1.50%  '

# Row 35
$ws.Range('D35').Value = $q + '3.33'
$ws.Range('E35').Value = $q + '  -1.79%  '

# Row 36
$ws.Range('D36').Value = $q + '19.00'
$ws.Range('E36').Value = $q + '  +15.28%  '

# Row 37
$ws.Range('D37').Value = $q + '26.28'
$ws.Range('E37').Value = $q + '  +12.68%  '

# Row 38
$ws.Range('E38').Value = $q + '  -0.99%  '

# Row 39
$ws.Range('E39').Value = $q + '  -0.31%  '

# Row 40
$ws.Range('B40').Value = $q + 'NEARProtocol'
$ws.Range('C40').Value = $q + 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').Value = $q + '3.48'
$ws.Range('E40').Value = $q + '  +0.54%  '

# Row 41
$ws.Range('B41').Value = $q + 'ApeXProtocol'
$ws.Range('C41').Value = $q + 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D41').Value = $q + '2.11'
$ws.Range('E41').Value = $q + '  +32.99%  '

# Row 42
$ws.Range('D42').Value = $q + '3.88'
$ws.Range('E42').Value = $q + '  +0.36%  '

# Row 43
$ws.Range('D43').Value = $q + '2.094.16'
$ws.Range('E43').Value = $q + '  +1.22%  '

# Row 44
$ws.Range('E44').Value = $q + '  -2.55%  '

# Row 45
$ws.Range('E45').Value = $q + '  +0.04%  '

# Row 46
$ws.Range('D46').Value = $q + '86.59'
$ws.Range('E46').Value = $q + '  +1.08%  '

# Row 47
$ws.Range('D47').Value = $q + '9.10'
$ws.Range('E47').Value = $q + '  +4.00%  '

# Row 48
$ws.Range('B48').Value = $q + 'ordi'
$ws.Range('C48').Value = $q + 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').Value = $q + '75.37'
$ws.Range('E48').Value = $q + '  +8.88%  '

# Row 49
$ws.Range('B49').Value = $q + 'RocketPoolETH'
$ws.Range('C49').Value = $q + 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = $q + '2.803.37'
$ws.Range('E49').Value = $q + '  +0.26%  '

# Row 50
$ws.Range('D50').Value = $q + '103.96'
$ws.Range('E50').Value = $q + '  -0.30%  '

# Row 51
$ws.Range('E51').Value = $q + '  +2.39%  '
